$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9693716918425304
$ws.Range("J2").Value = 0.9693716918425304
$ws.Range("M2").Value = 14.48297233333333
$ws.Range("N2").Value = 43.448917
$ws.Range("O2").Value = 0.2019336017030403
$ws.Range("P2").Value = 0.2019336017030403
$ws.Range("Q2").Value = 132.6199886821252
$ws.Range("R2").Value = 1193.579898139126
$ws.Range("S2").Value = 0.1957487171227319
$ws.Range("T2").Value = 0.1957487171227319
$ws.Range("I3").Value = 0.9693716918425304
$ws.Range("J3").Value = 0.9693716918425304
$ws.Range("O3").Value = 0.007144147385663391
$ws.Range("P3").Value = 0.00714414738566339
$ws.Range("S3").Value = 0.006925334238012911
$ws.Range("T3").Value = 0.006925334238012911
$ws.Range("I4").Value = 0.9693716918425304
$ws.Range("J4").Value = 0.9693716918425304
$ws.Range("M4").Value = 56.726097
$ws.Range("N4").Value = 170.178291
$ws.Range("O4").Value = 0.7909222509112964
$ws.Range("P4").Value = 0.7909222509112963
$ws.Range("Q4").Value = 519.4385633677221
$ws.Range("R4").Value = 4674.947070309498
$ws.Range("S4").Value = 0.7666976404817857
$ws.Range("T4").Value = 0.7666976404817856
$ws.Range("G5").Value = 0.2893236666666667
$ws.Range("H5").Value = 0.867971
$ws.Range("I5").Value = 0.03062830815746963
$ws.Range("J5").Value = 0.03062830815746962
$ws.Range("M5").Value = 14.48297233333333
$ws.Range("N5").Value = 43.448917
$ws.Range("O5").Value = 0.2019336017030403
$ws.Range("P5").Value = 0.2019336017030403
$ws.Range("Q5").Value = 4.19026665971189
$ws.Range("R5").Value = 37.71239993740701
$ws.Range("S5").Value = 0.006184884580308452
$ws.Range("T5").Value = 0.006184884580308451
$ws.Range("G6").Value = 0.2893236666666667
$ws.Range("H6").Value = 0.867971
$ws.Range("I6").Value = 0.03062830815746963
$ws.Range("J6").Value = 0.03062830815746962
$ws.Range("O6").Value = 0.007144147385663391
$ws.Range("P6").Value = 0.00714414738566339
$ws.Range("Q6").Value = 0.1482461677984445
$ws.Range("R6").Value = 1.334215510186
$ws.Range("S6").Value = 0.0002188131476504793
$ws.Range("T6").Value = 0.0002188131476504793
$ws.Range("G7").Value = 0.2893236666666667
$ws.Range("H7").Value = 0.867971
$ws.Range("I7").Value = 0.03062830815746963
$ws.Range("J7").Value = 0.03062830815746962
$ws.Range("M7").Value = 56.726097
$ws.Range("N7").Value = 170.178291
$ws.Range("O7").Value = 0.7909222509112964
$ws.Range("P7").Value = 0.7909222509112963
$ws.Range("Q7").Value = 16.412202379729
$ws.Range("R7").Value = 147.709821417561
$ws.Range("S7").Value = 0.0242246104295107
$ws.Range("T7").Value = 0.02422461042951069
